# Populate the "Level" column (C) with the review-feature difficulty levels.
# Row layout (unchanged): A=English, B=Tieng Viet, C=Level, D=PathImage
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$levels = @{
    2  = "easy"      # cluster / cụm
    3  = "hard"      # residential area / khu dân cư
    4  = "moderate"  # hamlet / thôn, ấp, xóm
    5  = "hard"      # commune / xã
    6  = "easy"       # ward / phường
    7  = "easy"       # building / tòa nhà
    8  = "easy"       # urban zone / khu vực đô thị
    9  = "easy"       # province / tỉnh
    10 = "easy"       # corner / góc ngã tư
    11 = "easy"       # street / phố
    12 = "nothing"    # grandfather / ông (nội, ngoại)
    13 = "nothing"    # avenue / đại lộ
    14 = "nothing"    # highway / đường cao tốc
    15 = "nothing"    # national route / quốc lộ
    16 = "nothing"    # sister / chị/em gái ruột
    17 = "nothing"    # lane / ngõ, hẻm; làn đường
    18 = "nothing"    # tower / tháp
    19 = "nothing"    # grandmother / bà
    20 = "nothing"    # ancestor / tổ tiên
    21 = "nothing"    # provincial route / tỉnh lộ
    22 = "nothing"    # industrial park / khu công nghiệp
    23 = "nothing"    # great-grandparent / ông bà cố
    24 = "nothing"    # district / quận
    25 = "nothing"    # family / gia đình
    26 = "nothing"    # brother / anh/em trai ruột
    27 = "nothing"    # quarter / khu phố, tổ
    28 = "nothing"    # export processing zone / khu chế xuất
}

foreach ($row in $levels.Keys) {
    $ws.Range("C$row").Value = $levels[$row]
}
